$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.543.05"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "1.942.19"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.46"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.38"
$ws.Range("E8").Value = "  -1.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.361"
$ws.Range("E9").Value = "  -2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0846"
$ws.Range("E10").Value = "  +0.39%  "

$ws.Range("E11").Value = "  -0.92%  "

$ws.Range("D12").Value = "2.229.59"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.30"
$ws.Range("E13").Value = "  -3.64%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.812"
$ws.Range("E14").Value = "  -2.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.43"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("E16").Value = "  -3.13%  "

$ws.Range("D17").Value = "1.933.88"
$ws.Range("E17").Value = "  -1.60%  "

$ws.Range("D18").Value = "36.465.31"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.35"
$ws.Range("E19").Value = "  -2.92%  "

$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  -2.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.06"
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.99"
$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("E24").Value = "  -6.61%  "

$ws.Range("E25").Value = "  +1.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.20"
$ws.Range("E26").Value = "  -3.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.49"
$ws.Range("E27").Value = "  -2.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("E28").Value = "  +8.25%  "

$ws.Range("E29").Value = "  -3.66%  "

$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("E31").Value = "  -4.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.59"
$ws.Range("E32").Value = "  -3.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0617"
$ws.Range("E33").Value = "  -3.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.19"
$ws.Range("E34").Value = "  -2.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.23"
$ws.Range("E35").Value = "  +4.62%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  -1.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.17"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.18"
$ws.Range("E39").Value = "  +9.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0985"
$ws.Range("E40").Value = "  +2.53%  "

$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.01"
$ws.Range("E44").Value = "  +1.67%  "

$ws.Range("D45").Value = "1.342.88"

$ws.Range("E46").Value = "  -2.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.64"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("E48").Value = "  -1.29%  "

$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").Value = "2.120.60"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.26"
$ws.Range("E51").Value = "  -3.68%  "
